$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5249140858650208
$ws.Range("B1").Value = 1.221789598464966
$ws.Range("C1").Value = 1.712547183036804
$ws.Range("D1").Value = -1
$ws.Range("E1").Value = -1
